# Add authentication error handling test row to the student workbook.
# Appends a new student-style record (row 13) used to exercise the
# "BSU / A1" authentication error-handling path added in this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "123 123"
$ws.Range("B13").Value = "BSU"
$ws.Range("C13").Value = "A1"
$ws.Range("D13").Value = 2023
